$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: EdgeXFoundry data ---
$ws.Range("D2").Value = "Linux Foundation (Edge)"

$ws.Range("D3").Value = "https://lfedge.org/projects/edgex-foundry/"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://lfedge.org/projects/edgex-foundry/")

$ws.Range("D4").Value = "https://github.com/edgexfoundry/edgex-go"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://github.com/edgexfoundry/edgex-go")

$ws.Range("D6").Value = 481
$ws.Range("D7").Value = "1.3k"
$ws.Range("D8").Value = "4.5k"
$ws.Range("D9").Value = "Apache-2.0"

$ws.Range("D5").Value = "https://docs.edgexfoundry.org/3.1/"
$ws.Hyperlinks.Add($ws.Range("D5"), "https://docs.edgexfoundry.org/3.1/")

$ws.Range("D11").Value = "https://docs.edgexfoundry.org/3.1/microservices/device/services/device-modbus/"

$ws.Range("D10").Value = "✓"
$ws.Hyperlinks.Add($ws.Range("D10"), "https://docs.edgexfoundry.org/3.1/microservices/device/services/device-modbus/")

Write-Output "done values"
